# Auto-generated from the unified diff of xl/worksheets/sheet1.xml.
# Updates the crypto price/volume table cells that changed between commits,
# including the two pairs of rows whose coin data was fully replaced
# (rows 35/36: ARBITRUM <-> LidoDAOToken, rows 44/45: RocketPoolETH <-> FraxShare).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, exactly as stored in the source XML
# (<is><t>...</t></is> inline strings), without letting Excel silently
# reinterpret numeric-looking strings (e.g. "0.9998", "1.000") as numbers.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.199.19"
Set-TextValue "E2" "  -0.12%  "
# Row 3
Set-TextValue "D3" "1.855.67"
Set-TextValue "E3" "  -0.26%  "
# Row 4
Set-TextValue "D4" "0.9998"
Set-TextValue "E4" "  +0.04%  "
# Row 5
Set-TextValue "D5" "241.35"
Set-TextValue "E5" "  -0.43%  "
# Row 6
Set-TextValue "D6" "0.6987"
Set-TextValue "E6" "  -0.88%  "
# Row 7
Set-TextValue "D7" "0.9999"
Set-TextValue "E7" "  +0.01%  "
# Row 8
Set-TextValue "D8" "0.07779"
Set-TextValue "E8" "  -0.31%  "
# Row 9
Set-TextValue "D9" "0.3076"
Set-TextValue "E9" "  -1.60%  "
# Row 10
Set-TextValue "D10" "23.72"
Set-TextValue "E10" "  -2.17%  "
# Row 11
Set-TextValue "D11" "0.07807"
Set-TextValue "E11" "  -2.39%  "
# Row 12
Set-TextValue "D12" "1.863.86"
Set-TextValue "E12" "  +0.16%  "
# Row 13
Set-TextValue "D13" "5.108"
Set-TextValue "E13" "  -1.39%  "
# Row 14
Set-TextValue "D14" "92.12"
Set-TextValue "E14" "  -2.13%  "
# Row 15
Set-TextValue "D15" "0.6874"
Set-TextValue "E15" "  -1.40%  "
# Row 16
Set-TextValue "D16" "6.544"
Set-TextValue "E16" "  +2.30%  "
# Row 17
Set-TextValue "D17" "0.000008471"
Set-TextValue "E17" "  +2.09%  "
# Row 18
Set-TextValue "D18" "29.203.75"
Set-TextValue "E18" "  -0.22%  "
# Row 19
Set-TextValue "D19" "247.98"
Set-TextValue "E19" "  -2.15%  "
# Row 20
Set-TextValue "D20" "2.107.18"
Set-TextValue "E20" "  -0.27%  "
# Row 21
Set-TextValue "D21" "12.83"
Set-TextValue "E21" "  -2.22%  "
# Row 22
Set-TextValue "D22" "1.000"
# Row 23
Set-TextValue "D23" "7.545"
Set-TextValue "E23" "  +0.24%  "
# Row 24
Set-TextValue "D24" "0.9999"
Set-TextValue "E24" "  +0.09%  "
# Row 25
Set-TextValue "D25" "0.1506"
Set-TextValue "E25" "  -3.53%  "
# Row 26
Set-TextValue "D26" "161.49"
Set-TextValue "E26" "  +1.16%  "
# Row 27
Set-TextValue "D27" "8.859"
Set-TextValue "E27" "  -1.48%  "
# Row 28
Set-TextValue "E28" "  -2.14%  "
# Row 29
Set-TextValue "D29" "1.556"
Set-TextValue "E29" "  +3.82%  "
# Row 30
Set-TextValue "D30" "4.257"
Set-TextValue "E30" "  -1.20%  "
# Row 31
Set-TextValue "D31" "4.209"
Set-TextValue "E31" "  -1.34%  "
# Row 32
Set-TextValue "E32" "  -1.33%  "
# Row 33
Set-TextValue "D33" "0.05233"
Set-TextValue "E33" "  -0.78%  "
# Row 34
Set-TextValue "D34" "0.7610"
Set-TextValue "E34" "  +1.73%  "
# Row 35
Set-TextValue "B35" "LidoDAOToken"
Set-TextValue "C35" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D35" "1.846"
Set-TextValue "E35" "  -2.10%  "
# Row 36
Set-TextValue "B36" "ARBITRUM"
Set-TextValue "C36" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D36" "1.169"
Set-TextValue "E36" "  +1.19%  "
# Row 37
Set-TextValue "D37" "2.709"
Set-TextValue "E37" "  -0.10%  "
# Row 38
Set-TextValue "E38" "  -0.37%  "
# Row 39
Set-TextValue "D39" "1.224.18"
Set-TextValue "E39" "  -1.71%  "
# Row 40
Set-TextValue "D40" "2.730"
Set-TextValue "E40" "  -0.15%  "
# Row 41
Set-TextValue "D41" "0.9010"
Set-TextValue "E41" "  +0.64%  "
# Row 42
Set-TextValue "D42" "109.12"
Set-TextValue "E42" "  -1.66%  "
# Row 43
Set-TextValue "D43" "0.9993"
Set-TextValue "E43" "  +0.03%  "
# Row 44
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "5.504"
Set-TextValue "E44" "  -10.31%  "
# Row 45
Set-TextValue "B45" "RocketPoolETH"
Set-TextValue "C45" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D45" "2.004.11"
Set-TextValue "E45" "  -0.09%  "
# Row 46
Set-TextValue "E46" "  -3.04%  "
# Row 47
Set-TextValue "D47" "65.41"
Set-TextValue "E47" "  -7.59%  "
# Row 48
Set-TextValue "D48" "0.5181"
Set-TextValue "E48" "  -0.09%  "
# Row 49
Set-TextValue "D49" "9.536"
Set-TextValue "E49" "  +0.61%  "
# Row 50
Set-TextValue "D50" "1.750"
Set-TextValue "E50" "  -1.96%  "
# Row 51
Set-TextValue "D51" "7.050"
Set-TextValue "E51" "  +0.36%  "
